$wb = $excel.ActiveWorkbook

# Update the Metadata sheet: URL and Date values
$wsMeta = $wb.Worksheets.Item("Metadata")
$wsMeta.Range("B2").Value = "http://example.org/ig/exampleig/ValueSet/gvhd-all-valueset"
$wsMeta.Range("B8").Value = "2023-04-27T11:07:01-05:00"

# Update the "Include ValueSets" sheet: ValueSet URL for SCT codes
$wsInclude1 = $wb.Worksheets.Item("Include ValueSets")
$wsInclude1.Range("A2").Value = "http://example.org/ig/exampleig/ValueSet/gvhd-sct-codes"

# Update the "Include ValueSets 2" sheet: ValueSet URL for ICD10 codes
$wsInclude2 = $wb.Worksheets.Item("Include ValueSets 2")
$wsInclude2.Range("A2").Value = "http://example.org/ig/exampleig/ValueSet/gvhd-icd10-codes"
